$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has a header row (row 1) and Nelson Ouya's attendance
# in row 2. We need to insert Pelegrin Ogalo's attendance as the new row 2,
# pushing Nelson Ouya's existing row down to row 3.

# 1) Copy Nelson Ouya's existing row 2 values down to row 3 (use Value2 to
#    read back the real underlying value).
$ws.Range("A3").Value = $ws.Range("A2").Value2
$ws.Range("B3").Value = $ws.Range("B2").Value2
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = $ws.Range("C2").Value2
$ws.Range("C3").ClearFormats()
$ws.Range("D3").Value = $ws.Range("D2").Value2
$ws.Range("E3").Value = $ws.Range("E2").Value2
$ws.Range("F3").Value = $ws.Range("F2").Value2

# 2) Overwrite row 2 with Pelegrin Ogalo's attendance for 2025-05-06.
$ws.Range("A2").Value = "Pelegrin Ogalo"
$ws.Range("B2").Value = "pelegrin@gmail.com"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "2025-05-06"
$ws.Range("C2").ClearFormats()
$ws.Range("D2").Value = "Present"
$ws.Range("E2").Value = "N/A"
$ws.Range("F2").Value = "N/A"
